# Generate Report for Handoff
# Replaces the two stale "e2e/<guid>.md" + "<guid>.png" source entries with
# fresh entries for a .md handoff ("c12cd017-...") and its dependency
# ("dd08a102-..."), updates the handoff timestamps, and drops the old third
# (dependency/.png) row on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md"
$ov.Range("D2").Value = "2016-03-13 07:03:29"

$ov.Range("A3").Value = "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md"
$ov.Range("D3").Value = "2016-03-13 07:03:29"

$ov.Rows.Item(4).Delete()

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.md", [Type]::Missing, [Type]::Missing, "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/67bcde2d-49e4-41a1-b073-b278d84bb501.png", [Type]::Missing, [Type]::Missing, "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md"
$zh.Range("D2").Value = "c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-13 07:03:26"

$zh.Range("A3").Value = "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md"
$zh.Range("B3").Value = ".md"
$zh.Range("D3").Value = "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-13 07:03:26"
$zh.Range("I3").Value = "Include"
$zh.Range("J3").Value = ""

$zh.Rows.Item(4).Delete()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.md", [Type]::Missing, [Type]::Missing, "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2593e9b4dd3ff381ac1d4336114c65889270cc57/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.f0c2e2c83004b4109dcf385861a1520489feb02c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/67bcde2d-49e4-41a1-b073-b278d84bb501.png", [Type]::Missing, [Type]::Missing, "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/67bcde2d-49e4-41a1-b073-b278d84bb501.png", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2593e9b4dd3ff381ac1d4336114c65889270cc57/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5e17b051bc8da2614b90dc6b43aa50b1a087c274.png", [Type]::Missing, [Type]::Missing, "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md"
$de.Range("D2").Value = "c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.de-de.xlf"
$de.Range("E2").Value = "2016-03-13 07:03:29"

$de.Range("A3").Value = "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md"
$de.Range("B3").Value = ".md"
$de.Range("D3").Value = "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.de-de.xlf"
$de.Range("E3").Value = "2016-03-13 07:03:29"
$de.Range("I3").Value = "Include"
$de.Range("J3").Value = ""

$de.Rows.Item(4).Delete()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.md", [Type]::Missing, [Type]::Missing, "c12cd017-6b27-4c1a-a737-cf9612fd3db3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e458c96d555491aa2eadb555907f4dfdb0f3910/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2d810e9b-ad4c-4742-a96f-db848a5fc7e5.f0c2e2c83004b4109dcf385861a1520489feb02c.de-de.xlf", [Type]::Missing, [Type]::Missing, "c12cd017-6b27-4c1a-a737-cf9612fd3db3.d4dcba8639963b2215a27b38f8a847e51f789549.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/67bcde2d-49e4-41a1-b073-b278d84bb501.png", [Type]::Missing, [Type]::Missing, "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/e7806561b5c61dfae923047844b16820ea409292/e2e/67bcde2d-49e4-41a1-b073-b278d84bb501.png", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e458c96d555491aa2eadb555907f4dfdb0f3910/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5e17b051bc8da2614b90dc6b43aa50b1a087c274.png", [Type]::Missing, [Type]::Missing, "dd08a102-2d58-4b5b-bcf7-e8440b9b0360.ce40297b13a84debd09da7a3288f9090003c73ad.de-de.xlf") | Out-Null

$wb.Save()
